# Auto-generated edit script: updates Faerie Profits market data values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1001.95
$ws.Range("I28").Value = 1147.6875
$ws.Range("J28").Value = 419
$ws.Range("K28").Value = 1147.6875
$ws.Range("L28").Value = 419
$ws.Range("M28").Value = -662.6875
$ws.Range("N28").Value = -1389

# Row 80
$ws.Range("H80").Value = 1672.5641
$ws.Range("J80").Value = 1949.5652
$ws.Range("L80").Value = 5848.6956
$ws.Range("N80").Value = -7844.6956

# Row 83
$ws.Range("H83").Value = 1672.5641
$ws.Range("J83").Value = 1949.5652
$ws.Range("L83").Value = 17546.0868
$ws.Range("N83").Value = -27530.0868

# Row 86
$ws.Range("H86").Value = 4074.4075
$ws.Range("I86").Value = 1414.1428
$ws.Range("K86").Value = 1414.1428
$ws.Range("M86").Value = -291.1428000000001

# Row 89
$ws.Range("H89").Value = 4074.4075
$ws.Range("I89").Value = 1414.1428
$ws.Range("K89").Value = 7070.714
$ws.Range("M89").Value = -1454.714

# Row 132
$ws.Range("H132").Value = 35718170
$ws.Range("I132").Value = 37040876
$ws.Range("K132").Value = 111122628
$ws.Range("M132").Value = -111120098

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 931.2286
$ws.Range("I2").Value = 964.7917
$ws.Range("J2").Value = 858
$ws.Range("K2").Value = 964.7917
$ws.Range("L2").Value = 858
$ws.Range("M2").Value = -851.7917
$ws.Range("N2").Value = -1084

# Row 32
$ws.Range("H32").Value = 1887.06
$ws.Range("I32").Value = 1887.06
$ws.Range("K32").Value = 1887.06
$ws.Range("M32").Value = -1600.06

# Row 61
$ws.Range("H61").Value = 884887.5
$ws.Range("I61").Value = 1321370.9
$ws.Range("K61").Value = 1321370.9
$ws.Range("M61").Value = -1321158.9

# Row 110
$ws.Range("H110").Value = 1831.25
$ws.Range("I110").Value = 1831.25
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1831.25
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 213.75
$ws.Range("N110").ClearContents()

# Row 116
$ws.Range("H116").Value = 931.2286
$ws.Range("I116").Value = 964.7917
$ws.Range("J116").Value = 858
$ws.Range("K116").Value = 964.7917
$ws.Range("L116").Value = 858
$ws.Range("M116").Value = 1329.2083
$ws.Range("N116").Value = -5446

# Row 132
$ws.Range("H132").Value = 2851151.8
$ws.Range("I132").Value = 3346534.5
$ws.Range("K132").Value = 10039603.5
$ws.Range("M132").Value = -10037073.5

# Row 136
$ws.Range("H136").Value = 884887.5
$ws.Range("I136").Value = 1321370.9
$ws.Range("K136").Value = 3964112.7
$ws.Range("M136").Value = -3961562.7

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 931.2286
$ws.Range("I3").Value = 964.7917
$ws.Range("J3").Value = 858
$ws.Range("K3").Value = 964.7917
$ws.Range("L3").Value = 858
$ws.Range("M3").Value = -850.7917
$ws.Range("N3").Value = -1086

# Row 81
$ws.Range("H81").Value = 39686.668
$ws.Range("J81").Value = 39686.668
$ws.Range("L81").Value = 39686.668
$ws.Range("N81").Value = -41808.668

# Row 84
$ws.Range("H84").Value = 39686.668
$ws.Range("J84").Value = 39686.668
$ws.Range("L84").Value = 119060.004
$ws.Range("N84").Value = -129668.004

# Row 86
$ws.Range("H86").Value = 2473490.5
$ws.Range("I86").Value = 3924374.5
$ws.Range("K86").Value = 3924374.5
$ws.Range("M86").Value = -3923251.5

# Row 89
$ws.Range("H89").Value = 2473490.5
$ws.Range("I89").Value = 3924374.5
$ws.Range("K89").Value = 19621872.5
$ws.Range("M89").Value = -19616256.5

# Row 99
$ws.Range("H99").Value = 3894.5
$ws.Range("J99").Value = 4932.6665
$ws.Range("L99").Value = 4932.6665
$ws.Range("N99").Value = -7928.6665

# Row 107
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

# Row 134
$ws.Range("H134").Value = 8136.102
$ws.Range("I134").Value = 4419.706
$ws.Range("K134").Value = 13259.118
$ws.Range("M134").Value = -10724.118

# Row 135
$ws.Range("H135").Value = 102260
$ws.Range("J135").Value = 102260
$ws.Range("L135").Value = 102260
$ws.Range("N135").Value = -112400

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2274.8
$ws.Range("I31").Value = 2274.8
$ws.Range("K31").Value = 2274.8
$ws.Range("M31").Value = -1979.8

# Row 34
$ws.Range("H34").Value = 2274.8
$ws.Range("I34").Value = 2274.8
$ws.Range("K34").Value = 2274.8
$ws.Range("M34").Value = -2072.8

# Row 122
$ws.Range("H122").Value = 4270.826
$ws.Range("I122").Value = 4182.8945
$ws.Range("J122").Value = 4688.5
$ws.Range("K122").Value = 12548.6835
$ws.Range("L122").Value = 14065.5
$ws.Range("M122").Value = -10098.6835
$ws.Range("N122").Value = -18965.5

# Row 132
$ws.Range("H132").Value = 1145977.4
$ws.Range("I132").Value = 2224665
$ws.Range("J132").Value = 3837.647
$ws.Range("K132").Value = 6673995
$ws.Range("L132").Value = 11512.941
$ws.Range("M132").Value = -6671465
$ws.Range("N132").Value = -16572.941

$ws = $wb.Worksheets.Item("CUL")
# Row 28
$ws.Range("H28").Value = 2015.5
$ws.Range("I28").Value = 698
$ws.Range("K28").Value = 2094
$ws.Range("M28").Value = -1862

# Row 127
$ws.Range("H127").Value = 3297
$ws.Range("J127").Value = 3297
$ws.Range("L127").Value = 9891
$ws.Range("N127").Value = -19811

# Row 129
$ws.Range("H129").Value = 63630.562
$ws.Range("I129").Value = 83921.414
$ws.Range("J129").Value = 2758
$ws.Range("K129").Value = 251764.242
$ws.Range("L129").Value = 8274
$ws.Range("M129").Value = -246764.242
$ws.Range("N129").Value = -18274

# Row 131
$ws.Range("H131").Value = 190159.92
$ws.Range("I131").Value = 1111632.1
$ws.Range("J131").Value = 1676.9773
$ws.Range("K131").Value = 3334896.3
$ws.Range("L131").Value = 5030.9319
$ws.Range("M131").Value = -3329856.3
$ws.Range("N131").Value = -15110.9319

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3564.318
$ws.Range("I80").Value = 3165.5557
$ws.Range("J80").Value = 3840.3845
$ws.Range("K80").Value = 3165.5557
$ws.Range("L80").Value = 3840.3845
$ws.Range("M80").Value = -2167.5557
$ws.Range("N80").Value = -5836.3845

# Row 83
$ws.Range("H83").Value = 3564.318
$ws.Range("I83").Value = 3165.5557
$ws.Range("J83").Value = 3840.3845
$ws.Range("K83").Value = 15827.7785
$ws.Range("L83").Value = 19201.9225
$ws.Range("M83").Value = -10835.7785
$ws.Range("N83").Value = -29185.9225

# Row 102
$ws.Range("H102").Value = 19072
$ws.Range("I102").Value = 3332.318
$ws.Range("K102").Value = 3332.318
$ws.Range("M102").Value = -1710.318

# Row 107
$ws.Range("H107").Value = 646.2174
$ws.Range("I107").Value = 534.55554
$ws.Range("K107").Value = 534.55554
$ws.Range("M107").Value = 1385.44446

# Row 122
$ws.Range("H122").Value = 1690.1765
$ws.Range("I122").Value = 1117.0834
$ws.Range("K122").Value = 3351.2502
$ws.Range("M122").Value = -901.2501999999999

# Row 123
$ws.Range("H123").Value = 44149.832
$ws.Range("J123").Value = 44149.832
$ws.Range("L123").Value = 44149.832
$ws.Range("N123").Value = -49049.832

# Row 132
$ws.Range("H132").Value = 5643.231
$ws.Range("I132").Value = 5799.364
$ws.Range("J132").Value = 4784.5
$ws.Range("K132").Value = 17398.092
$ws.Range("L132").Value = 14353.5
$ws.Range("M132").Value = -14868.092
$ws.Range("N132").Value = -19413.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2475.9473
$ws.Range("I22").Value = 516
$ws.Range("K22").Value = 516
$ws.Range("M22").Value = -221

# Row 27
$ws.Range("H27").Value = 2475.9473
$ws.Range("I27").Value = 516
$ws.Range("K27").Value = 516
$ws.Range("M27").Value = -409

# Row 46
$ws.Range("H46").Value = 7650.3335
$ws.Range("J46").Value = 9636.214
$ws.Range("L46").Value = 9636.214
$ws.Range("N46").Value = -10012.214

# Row 68
$ws.Range("H68").Value = 3675
$ws.Range("I68").Value = 3662.5
$ws.Range("K68").Value = 3662.5
$ws.Range("M68").Value = -2913.5

# Row 71
$ws.Range("H71").Value = 3675
$ws.Range("I71").Value = 3662.5
$ws.Range("K71").Value = 18312.5
$ws.Range("M71").Value = -14568.5

# Row 94
$ws.Range("H94").Value = 41557.5
$ws.Range("J94").Value = 41557.5
$ws.Range("L94").Value = 41557.5
$ws.Range("N94").Value = -42909.5

# Row 100
$ws.Range("H100").Value = 3689.0588
$ws.Range("I100").Value = 3246.818
$ws.Range("K100").Value = 3246.818
$ws.Range("M100").Value = -2705.818

# Row 122
$ws.Range("H122").Value = 7095.3076
$ws.Range("I122").Value = 6435.625
$ws.Range("K122").Value = 19306.875
$ws.Range("M122").Value = -16856.875

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 33554.15
$ws.Range("I81").Value = 86568.414
$ws.Range("J81").Value = 4637.273
$ws.Range("K81").Value = 173136.828
$ws.Range("L81").Value = 9274.546
$ws.Range("M81").Value = -172075.828
$ws.Range("N81").Value = -11396.546

# Row 84
$ws.Range("H84").Value = 33554.15
$ws.Range("I84").Value = 86568.414
$ws.Range("J84").Value = 4637.273
$ws.Range("K84").Value = 865684.14
$ws.Range("L84").Value = 46372.73
$ws.Range("M84").Value = -860380.14
$ws.Range("N84").Value = -56980.73

# Row 96
$ws.Range("H96").Value = 2744.9092
$ws.Range("I96").Value = 2199.8572
$ws.Range("K96").Value = 2199.8572
$ws.Range("M96").Value = -826.8571999999999

# Row 100
$ws.Range("H100").Value = 1070.9333
$ws.Range("I100").Value = 908.88464
$ws.Range("K100").Value = 1817.76928
$ws.Range("M100").Value = -1276.76928

# Row 132
$ws.Range("H132").Value = 1907.5834
$ws.Range("I132").Value = 1998.6666
$ws.Range("K132").Value = 5995.9998
$ws.Range("M132").Value = -3465.9998

